{"js": "const replacements = [\n  [\"923\u00f72=461, 1\", \"807\u00f78=100, 7\"],\n  [\"183\u00f74=45, 3\", \"171\u00f72=85, 1\"],\n  [\"127\u00f77=18, 1\", \"651\u00f72=325, 1\"],\n  [\"250\u00f79=27, 7\", \"701\u00f78=87, 5\"],\n  [\"442\u00f73=147, 1\", \"250\u00f78=31, 2\"],\n  [\"757\u00f74=189, 1\", \"759\u00f78=94, 7\"],\n  [\"303\u00f72=151, 1\", \"822\u00f73=274, 0\"],\n  [\"310\u00f76=51, 4\", \"267\u00f79=29, 6\"],\n  [\"878\u00f79=97, 5\", \"265\u00f79=29, 4\"],\n  [\"834\u00f72=417, 0\", \"116\u00f75=23, 1\"],\n  [\"244\u00f72=122, 0\", \"887\u00f75=177, 2\"],\n  [\"642\u00f72=321, 0\", \"371\u00f77=53, 0\"],\n  [\"895\u00f76=149, 1\", \"483\u00f74=120, 3\"],\n  [\"159\u00f77=22, 5\", \"969\u00f79=107, 6\"],\n  [\"482\u00f72=241, 0\", \"181\u00f75=36, 1\"],\n  [\"709\u00f72=354, 1\", \"938\u00f78=117, 2\"],\n  [\"702\u00f76=117, 0\", \"828\u00f74=207, 0\"],\n  [\"280\u00f73=93, 1\", \"539\u00f79=59, 8\"],\n  [\"852\u00f75=170, 2\", \"583\u00f79=64, 7\"],\n  [\"264\u00f72=132, 0\", \"195\u00f72=97, 1\"],\n  [\"310\u00f78=38, 6\", \"701\u00f78=87, 5\"],\n  [\"477\u00f76=79, 3\", \"114\u00f74=28, 2\"],\n  [\"479\u00f78=59, 7\", \"267\u00f76=44, 3\"],\n  [\"731\u00f77=104, 3\", \"209\u00f72=104, 1\"],\n  [\"519\u00f75=103, 4\", \"889\u00f72=444, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('923\u00f72=461, 1', '807\u00f78=100, 7'),\n    @('183\u00f74=45, 3', '171\u00f72=85, 1'),\n    @('127\u00f77=18, 1', '651\u00f72=325, 1'),\n    @('250\u00f79=27, 7', '701\u00f78=87, 5'),\n    @('442\u00f73=147, 1', '250\u00f78=31, 2'),\n    @('757\u00f74=189, 1', '759\u00f78=94, 7'),\n    @('303\u00f72=151, 1', '822\u00f73=274, 0'),\n    @('310\u00f76=51, 4', '267\u00f79=29, 6'),\n    @('878\u00f79=97, 5', '265\u00f79=29, 4'),\n    @('834\u00f72=417, 0', '116\u00f75=23, 1'),\n    @('244\u00f72=122, 0', '887\u00f75=177, 2'),\n    @('642\u00f72=321, 0', '371\u00f77=53, 0'),\n    @('895\u00f76=149, 1', '483\u00f74=120, 3'),\n    @('159\u00f77=22, 5', '969\u00f79=107, 6'),\n    @('482\u00f72=241, 0', '181\u00f75=36, 1'),\n    @('709\u00f72=354, 1', '938\u00f78=117, 2'),\n    @('702\u00f76=117, 0', '828\u00f74=207, 0'),\n    @('280\u00f73=93, 1', '539\u00f79=59, 8'),\n    @('852\u00f75=170, 2', '583\u00f79=64, 7'),\n    @('264\u00f72=132, 0', '195\u00f72=97, 1'),\n    @('310\u00f78=38, 6', '701\u00f78=87, 5'),\n    @('477\u00f76=79, 3', '114\u00f74=28, 2'),\n    @('479\u00f78=59, 7', '267\u00f76=44, 3'),\n    @('731\u00f77=104, 3', '209\u00f72=104, 1'),\n    @('519\u00f75=103, 4', '889\u00f72=444, 1'),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n\n"}
